$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<othe>"
$ws.Range("C2").Value = 53

# Row 3
$ws.Range("B3").Value = "<ethe>"
$ws.Range("C3").Value = 52

# Row 4
$ws.Range("C4").Value = 51

# Row 5
$ws.Range("C5").Value = 50

# Row 6
$ws.Range("C6").Value = 51

# Row 7
$ws.Range("C7").Value = 45

# Row 8
$ws.Range("B8").Value = "<the>"
$ws.Range("C8").Value = 51

# Row 9
$ws.Range("B9").Value = "<willie>"
$ws.Range("C9").Value = 49

# Row 10
$ws.Range("C10").Value = 49

# Row 11
$ws.Range("C11").Value = 52

# Row 12
$ws.Range("B12").Value = "<bouh>"
$ws.Range("C12").Value = 52

# Row 13
$ws.Range("B13").Value = "<forte>"
$ws.Range("C13").Value = 44

# Row 14
$ws.Range("C14").Value = 49

# Row 15
$ws.Range("B15").Value = "<ale>"
$ws.Range("C15").Value = 51

# Row 16
$ws.Range("B16").Value = "<long>"
$ws.Range("C16").Value = 51

# Row 17
$ws.Range("B17").Value = "<canme>"
$ws.Range("C17").Value = 56

# Row 18
$ws.Range("C18").Value = 45
